$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Emily Andrew Young's experience (row 6, column G) from 3 to 9 years
$ws.Range("G6").Value = 9

# Append a new record for Logan Bruce Scott (row 7)
$ws.Range("A7").Value = "Logan"
$ws.Range("B7").Value = "Bruce"
$ws.Range("C7").Value = "Scott"
$ws.Range("D7").Value = "logan.s@email.com"
$ws.Range("E7").Value = "105 Magnolia Ave, KY"
$ws.Range("F7").Value = "Mechanic"
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = "Special"

# Append a new record for Lee Bruce Scott (row 8), with a hyperlinked email address
$ws.Range("A8").Value = "Lee"
$ws.Range("B8").Value = "Bruce"
$ws.Range("C8").Value = "Scott"
$ws.Range("D8").Value = "lee.s@email.com"
$ws.Range("E8").Value = "105 Magnolia Ave, KY"
$ws.Range("F8").Value = "Mechanic"
$ws.Range("G8").Value = 11
$ws.Range("H8").Value = "Special"

$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:lee.s@email.com")

[void]$ws.Range("D8").Select()
